$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: B1 label changes from "d13N" to "d13C" (sources combined / corrected signature column)
$ws.Range("B1").Value = "d13C"

# Updated / corrected d13C values for all data rows (2-53)
$ws.Range("B2").Value = -30.901849585646296
$ws.Range("B3").Value = -30.785648109436927
$ws.Range("B4").Value = -30.5038345294879
$ws.Range("B5").Value = -30.097176827060352
$ws.Range("B6").Value = -29.884875777540838
$ws.Range("B7").Value = -30.642374802662005
$ws.Range("B8").Value = -30.962949931992885
$ws.Range("B9").Value = -30.591655745345264
$ws.Range("B10").Value = -30.633792580221638
$ws.Range("B11").Value = -30.582603648575983
$ws.Range("B12").Value = -30.046714786577475
$ws.Range("B13").Value = -30.32630727175073
$ws.Range("B14").Value = -30.792331327251823
$ws.Range("B15").Value = -30.916260679954547
$ws.Range("B16").Value = -30.348997880890277
$ws.Range("B17").Value = -30.050189318037663
$ws.Range("B18").Value = -29.97114520212939
$ws.Range("B19").Value = -30.02560534729258
$ws.Range("B20").Value = -30.24468282186304
$ws.Range("B21").Value = -29.80029299263284
$ws.Range("B22").Value = -31.294687529559837
$ws.Range("B23").Value = -31.412815105535344
$ws.Range("B24").Value = -31.47964618495289
$ws.Range("B25").Value = -31.246585599982655
$ws.Range("B26").Value = -31.83257871594934
$ws.Range("B27").Value = -29.913014551693458
$ws.Range("B28").Value = -30.850227884094373
$ws.Range("B29").Value = -31.503751406597683
$ws.Range("B30").Value = -31.452392077233277
$ws.Range("B31").Value = -31.526191204204604
$ws.Range("B32").Value = -31.38851365950103
$ws.Range("B33").Value = -31.505393261731264
$ws.Range("B34").Value = -30.790076258516198
$ws.Range("B35").Value = -30.91636879160518
$ws.Range("B36").Value = -30.448443263896458
$ws.Range("B37").Value = -31.407364908381176
$ws.Range("B38").Value = -31.09206426062053
$ws.Range("B39").Value = -30.437175788588597
$ws.Range("B40").Value = -30.7248836617061
$ws.Range("B41").Value = -30.6898029798869
$ws.Range("B42").Value = -30.059327853221216
$ws.Range("B43").Value = -30.668071198638835
$ws.Range("B44").Value = -30.798659795402045
$ws.Range("B45").Value = -30.787611109446622
$ws.Range("B46").Value = -30.083324307283526
$ws.Range("B47").Value = -30.634215805829825
$ws.Range("B48").Value = -30.074407755116685
$ws.Range("B49").Value = -30.16748489507606
$ws.Range("B50").Value = -30.842293734897407
$ws.Range("B51").Value = -30.83926086978785
$ws.Range("B52").Value = -30.93791561306908
$ws.Range("B53").Value = -30.52544128107008

# Formatting corrections: remove now-unneeded emphasis/border formatting
# A2:A3 lose the (invisible) border-only style, reverting to the plain default style used elsewhere
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# B4:C4 lose the bold red font, matching the normal style used by the rest of the data rows
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# B53:C53 revert fully to the default (unstyled) cell format
$ws.Range("A53").Copy()
$ws.Range("B53:C53").PasteSpecial(-4122)

# Selection moved to the header row (A4:XFD4) as left by the editing session
[void]$ws.Range("A4:XFD4").Select()
